# Updated symbol list -- refresh Price (D) and Volume(1h) (E) columns for the
# cryptocurrencies that moved, matching the upstream scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of the cells below hold plain-text values in the source sheet (e.g. "297.49",
# "3.10%") rather than real numbers, so each one is touched individually: format as
# Text first so the literal string is preserved instead of Excel silently parsing it
# into a number/percentage, write the new value, then drop the style back to the
# workbook default so no stray formatting is left behind.
$cellRefs = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "D24",
    "E24",
    "D25",
    "E25",
    "E26",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)

foreach ($addr in $cellRefs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "297.53"
$ws.Range("E2").Value = "3.06%"
$ws.Range("D3").Value = "41.33"
$ws.Range("E3").Value = "2.33%"
$ws.Range("D4").Value = "5.007"
$ws.Range("E4").Value = "-0.49%"
$ws.Range("D5").Value = "0.07539"
$ws.Range("E5").Value = "3.37%"
$ws.Range("D6").Value = "1.572"
$ws.Range("E6").Value = "3.50%"
$ws.Range("D7").Value = "0.9304"
$ws.Range("E7").Value = "1.34%"
$ws.Range("D8").Value = "2.401"
$ws.Range("E8").Value = "0.17%"
$ws.Range("D9").Value = "0.1209"
$ws.Range("E9").Value = "2.26%"
$ws.Range("D10").Value = "0.1828"
$ws.Range("E10").Value = "6.05%"
$ws.Range("D11").Value = "0.08856"
$ws.Range("E11").Value = "2.40%"
$ws.Range("D12").Value = "0.04076"
$ws.Range("E12").Value = "-2.39%"
$ws.Range("D13").Value = "0.1054"
$ws.Range("E13").Value = "0.04%"
$ws.Range("D14").Value = "0.001274"
$ws.Range("E14").Value = "0.38%"
$ws.Range("D15").Value = "0.005917"
$ws.Range("E15").Value = "1.41%"
$ws.Range("D16").Value = "3.343"
$ws.Range("E16").Value = "-1.67%"
$ws.Range("E17").Value = "1.86%"
$ws.Range("D18").Value = "0.3332"
$ws.Range("E18").Value = "1.30%"
$ws.Range("D19").Value = "7.957"
$ws.Range("E19").Value = "1.12%"
$ws.Range("D20").Value = "0.1417"
$ws.Range("E20").Value = "5.56%"
$ws.Range("D21").Value = "0.2961"
$ws.Range("E21").Value = "2.71%"
$ws.Range("D22").Value = "0.04053"
$ws.Range("E22").Value = "4.73%"
$ws.Range("D23").Value = "0.001264"
$ws.Range("D24").Value = "0.003902"
$ws.Range("E24").Value = "2.11%"
$ws.Range("D25").Value = "0.0001229"
$ws.Range("E25").Value = "-4.22%"
$ws.Range("E26").Value = "-0.09%"
$ws.Range("D38").Value = "0.02427"
$ws.Range("E38").Value = "4.67%"
$ws.Range("D39").Value = "0.05216"
$ws.Range("E39").Value = "4.98%"
$ws.Range("D40").Value = "0.005883"
$ws.Range("E40").Value = "-15.16%"
$ws.Range("D41").Value = "0.007785"
$ws.Range("E41").Value = "1.27%"
$ws.Range("E42").Value = "4.42%"
$ws.Range("D43").Value = "0.007346"
$ws.Range("E43").Value = "-0.06%"
$ws.Range("D44").Value = "0.007831"
$ws.Range("E44").Value = "10.82%"
$ws.Range("D45").Value = "0.2979"
$ws.Range("E45").Value = "-4.55%"
$ws.Range("D46").Value = "0.00006295"
$ws.Range("E46").Value = "-2.16%"
$ws.Range("E47").Value = "-0.39%"
$ws.Range("D48").Value = "0.04474"
$ws.Range("E48").Value = "424.04%"
$ws.Range("D49").Value = "0.004192"
$ws.Range("E49").Value = "-0.22%"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").Value = "-0.39%"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").Value = "-0.39%"

foreach ($addr in $cellRefs) {
    $ws.Range($addr).Style = "Normal"
}
